$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 19:22"

# --- Row 4: Estados Unidos - refreshed case numbers ---
$ws.Range("B4").Value = 944805
$ws.Range("C4").Value = 19573
$ws.Range("E4").Value = 781002
$ws.Range("G4").Value = 1001
$ws.Range("H4").Value = 53194

# --- Row 23: Irlanda - refreshed case numbers ---
$ws.Range("B23").Value = 18561
$ws.Range("C23").Value = 377
$ws.Range("E23").Value = 8265
$ws.Range("G23").Value = 49
$ws.Range("H23").Value = 1063

# --- Rows 53/54: Egipto overtakes Sudafrica in the ranking ---
# Row 53 becomes Egipto with refreshed numbers
$ws.Range("A53").Value = "Egipto"
$ws.Range("B53").Value = 4319
$ws.Range("C53").Value = 227
$ws.Range("D53").Value = 1114
$ws.Range("E53").Value = 2898
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 13
$ws.Range("H53").Value = 307

# Row 54 becomes Sudafrica, keeping its previous numbers
$ws.Range("A54").Value = "Sudafrica"
$ws.Range("B54").Value = 4220
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 1473
$ws.Range("E54").Value = 2668
$ws.Range("F54").Value = 36
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 79

# --- Row 56: Luxemburgo - refreshed case numbers ---
$ws.Range("D56").Value = 3088
$ws.Range("E56").Value = 538
$ws.Range("F56").Value = 26

# --- Row 70: Irak - refreshed case numbers ---
$ws.Range("B70").Value = 1763
$ws.Range("C70").Value = 55
$ws.Range("D70").Value = 1224
$ws.Range("E70").Value = 453

# --- Row 91: Republica de Chipre - refreshed case numbers ---
$ws.Range("D91").Value = 148
$ws.Range("E91").Value = 648

# --- Rows 116-119: Mali jumps ahead of Mayotte/Kenia/Mauricio ---
# Row 116 becomes Mali with refreshed numbers
$ws.Range("A116").Value = "Mali"
$ws.Range("B116").Value = 370
$ws.Range("C116").Value = 45
$ws.Range("D116").Value = 91
$ws.Range("E116").Value = 258
$ws.Range("F116").Value = 0
$ws.Range("H116").Value = 21

# Row 117 becomes Mayotte, keeping its previous numbers
$ws.Range("A117").Value = "Mayotte"
$ws.Range("B117").Value = 354
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 144
$ws.Range("E117").Value = 206
$ws.Range("F117").Value = 4
$ws.Range("H117").Value = 4

# Row 118 becomes Kenia, keeping its previous numbers
$ws.Range("A118").Value = "Kenia"
$ws.Range("B118").Value = 343
$ws.Range("C118").Value = 7
$ws.Range("D118").Value = 98
$ws.Range("E118").Value = 231
$ws.Range("F118").Value = 2
$ws.Range("H118").Value = 14

# Row 119 becomes Mauricio, keeping its previous numbers
$ws.Range("A119").Value = "Mauricio"
$ws.Range("B119").Value = 331
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 295
$ws.Range("E119").Value = 27
$ws.Range("F119").Value = 3
$ws.Range("H119").Value = 9
